$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 12 - this shifts existing rows 12-23 down to 13-24,
# matching the diff where the old row12..row23 data now lives in row13..row24 and
# a brand-new data row is inserted as the new row 12.
$ws.Rows.Item(12).Insert()

# Populate the newly inserted row 12 with its data.
$ws.Range("A12").Value2 = 7
$ws.Range("B12").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C12").Value = "Ñuble"
$ws.Range("D12").Value2 = 45240
$ws.Range("E12").Value2 = 16
$ws.Range("F12").Value = "Fruta"
$ws.Range("G12").Value2 = 100102
$ws.Range("H12").Value = "Cítricos"
$ws.Range("I12").Value2 = 100102006
$ws.Range("J12").Value = "Pomelo"
$ws.Range("K12").Value = "Start Ruby"
$ws.Range("L12").Value = "Primera"
$ws.Range("M12").Value2 = 50
$ws.Range("N12").Value2 = 16000
$ws.Range("O12").Value2 = 16000
$ws.Range("P12").Value2 = 16000
$ws.Range("Q12").Value = '$/caja 14 kilos granel'
$ws.Range("R12").Value = "Región de O'Higgins"
$ws.Range("S12").Value2 = 1143
$ws.Range("T12").Value2 = 14
